$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 135 — this shifts the existing rows
# 135..176 down to 136..177 (dates/prices intact) and grows the
# used range to A1:R177.
$ws.Rows.Item(135).Insert()

# Populate the newly inserted row 135 with the new weekly record
# (same market/category context as its neighbours, new date & prices).
$ws.Range("A135").Value = 8
$ws.Range("B135").Value = "Terminal La Palmera de La Serena"
$ws.Range("C135").Value = "Coquimbo"
$ws.Range("D135").Value = 44876
$ws.Range("E135").Value = 4
$ws.Range("F135").Value = 100112001
$ws.Range("G135").Value = "Berenjena"
$ws.Range("H135").Value = "Sin especificar"
$ws.Range("I135").Value = "Primera"
$ws.Range("J135").Value = 400
$ws.Range("K135").Value = 11000
$ws.Range("L135").Value = 12000
$ws.Range("M135").Value = 11500
$ws.Range("N135").Value = "$/caja 40 unidades"
$ws.Range("O135").Value = "Región de Arica y Parinacota"
$ws.Range("P135").Value = 288
$ws.Range("Q135").Value = 40
$ws.Range("R135").Value = "Hortaliza"
